# Apply corrections to sentences/prompts in the stimulus list (exp/stim/trials/stim_list.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('E6').Value = 'Mariano hablaba del tiempo'
$ws.Range('E10').Value = 'Manuela vendía el carro'
$ws.Range('E11').Value = 'Daniel iba a Bolivia.'
$ws.Range('E14').Value = 'Mariano hablaba del tiempo.'
$ws.Range('E18').Value = 'Manuela vendía el carro.'
$ws.Range('E22').Value = '¿Mariano hablaba del tiempo?'
$ws.Range('E26').Value = '¿Manuela vendía el carro?'
$ws.Range('E27').Value = '¿Por qué iba a Bolivia?'
$ws.Range('E28').Value = '¿Cuándo leía el libro?'
$ws.Range('E29').Value = '¿Por qué ama la navidad?'
$ws.Range('E30').Value = '¿Por qué hablaba del agua?'
$ws.Range('E31').Value = '¿Cuándo lleva el abrigo?'
$ws.Range('E32').Value = '¿Cuándo bebía el vino?'
$ws.Range('E33').Value = '¿Por qué abre el regalo?'
$ws.Range('E34').Value = '¿Cuándo vendía el carro?'
$ws.Range('E36').Value = 'La hermana lavaba el plato.'
$ws.Range('E37').Value = 'Mi madre llama al médico.'
$ws.Range('E40').Value = 'El bebé comía muy bien.'
$ws.Range('E41').Value = 'La amiga vive en Nueva York'
$ws.Range('E42').Value = 'Mi novio viene del lago.'
$ws.Range('E44').Value = 'La hermana lavaba el plato.'
$ws.Range('E45').Value = 'Mi madre llama al médico.'
$ws.Range('D48').Value = '¿Cómo comía el bebé? '
$ws.Range('E48').Value = 'El bebé comía muy bien.'
$ws.Range('E49').Value = 'La amiga vive en Nueva York'
$ws.Range('D50').Value = '¿De dónde viene tu novio?'
$ws.Range('E50').Value = 'Mi novio viene del lago.'
$ws.Range('E52').Value = '¿La hermana lavaba el plato?'
$ws.Range('E53').Value = '¿Mi madre llama al médico?'
$ws.Range('E56').Value = '¿El bebé comía muy bien?'
$ws.Range('E57').Value = '¿La amiga vive en Nueva York?'
$ws.Range('E58').Value = '¿Mi novio viene del lago?'
$ws.Range('E59').Value = '¿Cuándo mira la luna?'
$ws.Range('E60').Value = '¿Cuándo lavaba el plato?'
$ws.Range('E61').Value = '¿Cuándo llama al médico?'
$ws.Range('E62').Value = '¿Por qué oía el río?'
$ws.Range('E63').Value = '¿Por qué odia a la reina?'
$ws.Range('E64').Value = '¿Cuándo comía muy bien?'
$ws.Range('E65').Value = '¿Por qué vive en Nueva York?'
$ws.Range('E66').Value = '¿Por qué viene del lago?'
